$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug with empty/incorrect airline and destination values for Ostrava/Brno rows
$ws.Range("D208").Value = "Travel Service"
$ws.Range("D209").Value = "bmi regional"
$ws.Range("B210").Value = "MUNICH"
$ws.Range("D210").Value = "Lufthansa"
$ws.Range("B211").Value = "LONDON STANSTED"
$ws.Range("D211").Value = "Ryanair"
$ws.Range("D213").Value = "Tunisair"
$ws.Range("D217").Value = "Wizz Air"
$ws.Range("D221").Value = "Wizz Air"
$ws.Range("D223").Value = "Travel Service"
$ws.Range("D226").Value = "Tunisair"
$ws.Range("D227").Value = "Travel Service"
$ws.Range("D229").Value = "Tailwind Airlines"
$ws.Range("B230").Value = "DJERBA / OSTRAVA"
$ws.Range("D230").Value = "Tunisair"
$ws.Range("D231").Value = "Travel Service"
$ws.Range("D236").Value = "Tailwind Airlines"
$ws.Range("B238").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D238").Value = "CZECH AIRLINES (CSA)"
$ws.Range("B239").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D239").Value = "KLM ROYAL DUTCH AIRLINES"
$ws.Range("B240").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D240").Value = "DELTA AIR LINES"
$ws.Range("B241").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D241").Value = "KOREAN AIR"
$ws.Range("B242").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D242").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("D243").Value = "TAROM ROMANIAN AIRLINES"
$ws.Range("B244").Value = "Crete / Heraklion, N. Kazantzakis Apt. (HER)"
$ws.Range("B245").Value = "Košice"
$ws.Range("D245").Value = "CZECH AIRLINES (CSA)"
$ws.Range("B246").Value = "Košice"
$ws.Range("D246").Value = "KLM ROYAL DUTCH AIRLINES"
$ws.Range("B247").Value = "Košice"
$ws.Range("D247").Value = "DELTA AIR LINES"
$ws.Range("B248").Value = "Košice"
$ws.Range("D248").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B249").Value = "Košice"
$ws.Range("D249").Value = "TAROM ROMANIAN AIRLINES"
$ws.Range("B250").Value = "Milan / Bergamo, Milan Bergamo Airport"
$ws.Range("D250").Value = "RYANAIR"
$ws.Range("B251").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D251").Value = "CZECH AIRLINES (CSA)"
$ws.Range("B252").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D252").Value = "KLM ROYAL DUTCH AIRLINES"
$ws.Range("B253").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D253").Value = "KOREAN AIR"
$ws.Range("B254").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D254").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B255").Value = "Corfu / Kerkyra, I. Kapodistrias (CFU)"
$ws.Range("B256").Value = "London, Stansted (STN)"
$ws.Range("D256").Value = "RYANAIR"
$ws.Range("B257").Value = "Antalya, Antayla (AYT)"
$ws.Range("D257").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B258").Value = "ERCAN"
$ws.Range("D258").Value = "TAILWIND"
$ws.Range("B259").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D259").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B260").Value = "Zakynthos, Zakinthos Is (ZTH)"
$ws.Range("D260").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B261").Value = "Rhodes, Diagoras Airport (RHO)"
$ws.Range("D261").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B262").Value = "Split"
$ws.Range("D262").Value = "CZECH AIRLINES (CSA)"
$ws.Range("B263").Value = "Hurghada, Hurghada (HRG)"
$ws.Range("D263").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B264").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D264").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B265").Value = "Mallorca, Palma de Mallorca (PMI)"
$ws.Range("D265").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B266").Value = "Marsa Alam, Marsa Alam (RMF)"
$ws.Range("D266").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B267").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D267").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B268").Value = "Varna, Varna Airport (VAR)"
$ws.Range("D269").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B270").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("B271").Value = "Kos, Kos Island International Airport (KGS)"
$ws.Range("D271").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B272").Value = "Varna, Varna Airport (VAR)"
$ws.Range("D272").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B273").Value = "MONASTIR"
$ws.Range("D273").Value = ""
$ws.Range("B274").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D274").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B275").Value = "BRATISLAVA"
$ws.Range("D275").Value = "TAILWIND"
$ws.Range("B276").Value = "Djerba"
$ws.Range("D276").Value = ""
$ws.Range("B277").Value = "KAVALA"
$ws.Range("D277").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B278").Value = "ALMERIA"
$ws.Range("D278").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B279").Value = "Podgorica"
$ws.Range("D279").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B280").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D280").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B281").Value = "Burgas, Burgas Airport (BOJ)"
$ws.Range("D281").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B282").Value = "Warsawa"
$ws.Range("B283").Value = "Warsawa"
$ws.Range("D283").Value = "KLM ROYAL DUTCH AIRLINES"
$ws.Range("B284").Value = "Warsawa"
$ws.Range("D284").Value = "DELTA AIR LINES"
$ws.Range("B285").Value = "Warsawa"
$ws.Range("D285").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B286").Value = "Warsawa"
$ws.Range("D286").Value = "TAROM ROMANIAN AIRLINES"
$ws.Range("B287").Value = "Prague, Václav Havel Airport Prague (PRG)"
$ws.Range("D287").Value = "TRAVEL SERVICE / SMARTWINGS"
$ws.Range("B288").Value = "Podgorica"
$ws.Range("D288").Value = "TRAVEL SERVICE / SMARTWINGS"

# Remove now-redundant duplicate rows at the end of the Ostrava block
$ws.Rows("289:292").Delete()

# Reset scroll position and selection as saved by the author
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F11").Select()
